$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 409, shifting existing rows 409..490 down to 410..491
$ws.Rows.Item(409).Insert()

# Populate the newly inserted row 409 with the new weekly record
$ws.Cells.Item(409, 1).Value = 4
$ws.Cells.Item(409, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(409, 3).Value = "Los Lagos"
$ws.Cells.Item(409, 4).Value = 44641
$ws.Cells.Item(409, 5).Value = 10
$ws.Cells.Item(409, 6).Value = 100112033
$ws.Cells.Item(409, 7).Value = "Lechuga"
$ws.Cells.Item(409, 8).Value = "Escarola"
$ws.Cells.Item(409, 9).Value = "Primera"
$ws.Cells.Item(409, 10).Value = 300
$ws.Cells.Item(409, 11).Value = 12000
$ws.Cells.Item(409, 12).Value = 14000
$ws.Cells.Item(409, 13).Value = 13000
$ws.Cells.Item(409, 14).Value = "$/caja 15 unidades"
$ws.Cells.Item(409, 15).Value = "Región del Maule"
$ws.Cells.Item(409, 16).Value = 867
$ws.Cells.Item(409, 17).Value = 15
$ws.Cells.Item(409, 18).Value = "Hortaliza"
